$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had 4 data rows (2-5): the first two (Sending cluster = "ECs")
# are obsolete under the refreshed TPM data and are dropped entirely; the
# remaining two rows (formerly 4-5, Sending cluster = "FAPs") move up to
# rows 2-3 and get refreshed metric values from the new TPM run.
$ws.Rows("2:3").Delete()

# Row 2 (originally row 4) - refreshed values
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("P2").Value = 0.3930722505137151
$ws.Range("S2").Value = 0.393072250513715
$ws.Range("T2").Value = 0.3930722505137151

# Row 3 (originally row 5) - refreshed values
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("P3").Value = 0.6069277494862849
$ws.Range("Q3").Value = 17.47057240755766
$ws.Range("S3").Value = 0.6069277494862849
$ws.Range("T3").Value = 0.6069277494862849
